$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from an existing header cell (E1) to the new header cells
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats

# Set the header text for the new columns
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Set boolean FALSE values for rows 2-5 in columns F, G, H
$ws.Range("F2:H5").Value = $false
